# Generate Report for Handoff
# Adds a new handoff record (0cfdf614-...) as row 3 to the Overview, zh-cn
# and de-de tables, mirroring the existing 7a280c5f-... row.

$wb = $excel.ActiveWorkbook

$newMd        = "0cfdf614-052d-4451-9e45-397471809b64ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdDisplay = "e2e\0cfdf614-052d-4451-9e45-397471809b64ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdTarget  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/325bd4f4b7d08eac2bfa8038695b1ac4a2e4e4f5/e2e/0cfdf614-052d-4451-9e45-397471809b64ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

$newZhCnXlf = "0cfdf614-052d-4451-9e45-397471809b64oooooooooooooooooooooooooooooooooooooooo.d226165824d97f9f480d494b923cfa0a0ca30520.zh-cn.xlf"
$newDeDeXlf = "0cfdf614-052d-4451-9e45-397471809b64oooooooooooooooooooooooooooooooooooooooo.d226165824d97f9f480d494b923cfa0a0ca30520.de-de.xlf"

$statusText   = "Ready for handoff"
$hoDate       = "2016-08-29 10:26:20"
$zhCnHoDate   = "2016-08-29 10:26:15"
$deDeHoDate   = "2016-08-29 10:26:20"
$handbackDate = "0001-01-01 00:00:00"

$hyperColor = 15570276  # BGR for FF6495ED, matches existing HyperLink font color
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Style-Hyperlink($range) {
    $range.Font.Color = $hyperColor
    $range.Font.Underline = 2
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newMd
$wsOverview.Range("B3").Value = $newMdDisplay
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = $hoDate
$wsOverview.Range("G3").NumberFormat = $dateFormat

Style-Hyperlink $wsOverview.Range("B3")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newMdTarget, "", "", $newMdDisplay) | Out-Null

$wsOverview.Range("E3:F3").ColumnWidth = 16.25

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $newMd
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $newZhCnXlf
$wsZhCn.Range("H3").Value = $zhCnHoDate
$wsZhCn.Range("H3").NumberFormat = $dateFormat
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = $handbackDate
$wsZhCn.Range("K3").NumberFormat = $dateFormat
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

Style-Hyperlink $wsZhCn.Range("A3")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newMdTarget, "", "", $newMd) | Out-Null

$wsZhCn.Range("C3").ColumnWidth = 16.25

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $newMd
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $newDeDeXlf
$wsDeDe.Range("H3").Value = $deDeHoDate
$wsDeDe.Range("H3").NumberFormat = $dateFormat
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = $handbackDate
$wsDeDe.Range("K3").NumberFormat = $dateFormat
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

Style-Hyperlink $wsDeDe.Range("A3")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newMdTarget, "", "", $newMd) | Out-Null

$wsDeDe.Range("C3").ColumnWidth = 16.25
